$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "37.849.42"
$ws.Range("E2").Value = "  +0.32%  "
Set-TextValue "D3" "2.084.94"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  +0.02%  "
Set-TextValue "D5" "234.09"
$ws.Range("E5").Value = "  +0.65%  "
Set-TextValue "D6" "0.627"
$ws.Range("E6").Value = "  +0.29%  "
Set-TextValue "D7" "58.98"
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("E8").Value = "  -0.01%  "
Set-TextValue "D9" "0.396"
$ws.Range("E9").Value = "  +2.21%  "
Set-TextValue "D10" "0.0791"
$ws.Range("E10").Value = "  +0.78%  "
Set-TextValue "D11" "0.109"
$ws.Range("E11").Value = "  +3.49%  "
Set-TextValue "D12" "2.389.90"
$ws.Range("E12").Value = "  +0.20%  "
Set-TextValue "D13" "14.77"
$ws.Range("E13").Value = "  +1.89%  "
Set-TextValue "D14" "21.29"
$ws.Range("E14").Value = "  +1.69%  "
Set-TextValue "D15" "0.780"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("E16").Value = "  +1.72%  "
Set-TextValue "D17" "2.064.59"
$ws.Range("E17").Value = "  -0.73%  "
Set-TextValue "D18" "37.761.00"
$ws.Range("E18").Value = "  +0.20%  "
Set-TextValue "D19" "6.18"
$ws.Range("E19").Value = "  +0.11%  "
Set-TextValue "D20" "71.79"
$ws.Range("E20").Value = "  +1.55%  "
Set-TextValue "D21" "0.0₃0847"
$ws.Range("E21").Value = "  +3.38%  "
Set-TextValue "D22" "228.46"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E24").Value = "  -0.44%  "
Set-TextValue "D25" "2.42"
$ws.Range("E25").Value = "  +1.60%  "
Set-TextValue "D26" "171.27"
$ws.Range("E26").Value = "  +0.61%  "
Set-TextValue "D27" "9.50"
$ws.Range("E27").Value = "  +6.72%  "
Set-TextValue "D28" "0.138"
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("E29").Value = "  -0.24%  "
Set-TextValue "D30" "19.58"
$ws.Range("E31").Value = "  +2.47%  "
Set-TextValue "D32" "4.77"
$ws.Range("E32").Value = "  +2.78%  "
Set-TextValue "D33" "0.0637"
$ws.Range("E33").Value = "  +2.03%  "
Set-TextValue "D34" "4.71"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("E35").Value = "  +0.29%  "
Set-TextValue "D36" "3.46"
$ws.Range("E36").Value = "  +2.17%  "
$ws.Range("E37").Value = "  -0.34%  "
Set-TextValue "D38" "0.999"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  +0.24%  "
Set-TextValue "D40" "0.0980"
$ws.Range("E40").Value = "  -0.70%  "
Set-TextValue "D41" "99.44"
$ws.Range("E41").Value = "  +0.53%  "
$ws.Range("E42").Value = "  +2.94%  "
Set-TextValue "D43" "17.18"
$ws.Range("E43").Value = "  +10.22%  "
$ws.Range("E44").Value = "  -0.81%  "
Set-TextValue "D45" "1.453.57"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("E46").Value = "  -0.57%  "
Set-TextValue "D47" "4.19"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("E48").Value = "  +1.67%  "
Set-TextValue "D49" "7.39"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("E50").Value = "  -0.72%  "
Set-TextValue "D51" "2.275.64"
$ws.Range("E51").Value = "  +0.26%  "
